$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (D3, D4): "state" -> "status"
$ws.Range("D3").Value = "status"
$ws.Range("D4").Value = "status"

# Data cells (D5:D18): "normal" -> "healthy"
$ws.Range("D5:D18").Value = "healthy"
